# Region XI_ELECTRIFICATION.xlsx update
# "Used most updated status accomplishmnet files as of may"
#
# Effect (see commit diff): two new columns are introduced after the
# existing "Unnamed: 43" column (AR):
#   - AR keeps its header ("Unnamed: 43") but its per-row values are
#     replaced: most rows go blank, a handful of rows get a new
#     "BBM ..." tag, and ENERGIZATION rows (2-21) keep blank too (their
#     tag is "ongrid" but that lands in the new AS column instead).
#   - AS is a brand new column (header "Unnamed: 44"); for rows 2-21 it
#     gets the literal "ongrid"; everywhere else it is blank.
#   - AT is a brand new column (header "Unnamed: 45") that receives the
#     OLD AR value for every data row (2-275).
#   - AU is what used to be AS: the "Status as of July 4, 2025" column,
#     complete with its dropdown data validation, now shifted two
#     columns to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers: AQ=43 AR=44 AS=45 AT=46 AU=47
$colAR = 44
$colAS = 45
$colAT = 46
$colAU = 47

# ---------------------------------------------------------------------
# 1) Header row (row 1)
#    AR1 ("Unnamed: 43") is untouched.
#    AS1 / AT1 are new header cells ("Unnamed: 44" / "Unnamed: 45"),
#    styled like the rest of the bold header row.
#    AU1 receives the old AS1 header text ("Status as of July 4, 2025").
# ---------------------------------------------------------------------

$oldAS1 = $ws.Cells.Item(1, $colAS).Value2

$ws.Cells.Item(1, $colAS).Value = "Unnamed: 44"
$ws.Cells.Item(1, $colAT).Value = "Unnamed: 45"
$ws.Cells.Item(1, $colAU).Value = $oldAS1

# Copy the bold/centered/bordered header formatting from AQ1 onto the
# two new header cells (AS1 and AT1) without disturbing their values.
$ws.Cells.Item(1, 43).Copy()
$ws.Range($ws.Cells.Item(1, $colAS), $ws.Cells.Item(1, $colAT)).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Data rows (2-275)
#    Every row's existing AR value moves two columns right, into AT.
#    AR itself is cleared, then re-populated only for the special rows
#    below ("ongrid" tag for ENERGIZATION rows 2-21, "BBM ..." tags for
#    a handful of later rows).
# ---------------------------------------------------------------------

$bbmTags = @{
    222 = "bbm 2023 ONGRID"
    224 = "BBM 2025 UPGRADE"
    233 = "BBM 2025 UPGRADE"
    234 = "BBM 2025 UPGRADE"
    248 = "BBM 2025 UPGRADE"
    249 = "BBM 2025 UPGRADE"
    251 = "BBM 2025 UPGRADE"
    257 = "BBM 2024 ONGRID"
    258 = "BBM 2024 ONGRID"
    259 = "BBM 2024 UPGRADE"
    260 = "BBM 2024 UPGRADE"
    261 = "BBM 2025 SOLAR"
    262 = "BBM 2024 UPGRADE"
    267 = "BBM 2025 ONGRID"
    268 = "BBM 2025 ONGRID"
    269 = "BBM 2024 ONGRID"
    270 = "BBM 2024 ONGRID"
    271 = "BBM 2024 UPGRADE"
    273 = "BBM 2024 UPGRADE"
    274 = "BBM 2024 UPGRADE"
}

for ($r = 2; $r -le 275; $r++) {
    $arCell = $ws.Cells.Item($r, $colAR)
    $oldVal = $arCell.Value2

    if ($oldVal -ne $null -and $oldVal -ne "") {
        $ws.Cells.Item($r, $colAT).Value = $oldVal
    }

    $arCell.ClearContents()

    if ($r -ge 2 -and $r -le 21) {
        $ws.Cells.Item($r, $colAS).Value = "ongrid"
    }

    if ($bbmTags.ContainsKey($r)) {
        $ws.Cells.Item($r, $colAR).Value = $bbmTags[$r]
    }
}

# ---------------------------------------------------------------------
# 3) Data validation: move the dropdown from AS2:AS275 to AU2:AU275
# ---------------------------------------------------------------------

$oldRange = $ws.Range("AS2:AS275")
$formula1 = $oldRange.Validation.Formula1
$oldRange.Validation.Delete()

$newRange = $ws.Range("AU2:AU275")
$newRange.Validation.Add(3, 1, 1, $formula1)
$newRange.Validation.IgnoreBlank = $true
$newRange.Validation.InCellDropdown = $true
$newRange.Validation.ShowInput = $false
$newRange.Validation.ShowError = $false
